$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 234 (pushes old rows 234:324 down to 235:325),
# matching a new weekly price entry for Arándano (blue) that sits right
# before the existing 2023-03-21 record.
$ws.Rows("234:234").Insert()

$ws.Cells.Item(234, 1).Value = 9
$ws.Cells.Item(234, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(234, 3).Value = "Metropolitana"
$ws.Cells.Item(234, 4).Value = 45007
$ws.Cells.Item(234, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(234, 5).Value = 13
$ws.Cells.Item(234, 6).Value = "Fruta"
$ws.Cells.Item(234, 7).Value = 100101
$ws.Cells.Item(234, 8).Value = "Berries"
$ws.Cells.Item(234, 9).Value = 100101001
$ws.Cells.Item(234, 10).Value = "Arándano (blue)"
$ws.Cells.Item(234, 11).Value = "Sin especificar"
$ws.Cells.Item(234, 12).Value = "Primera"
$ws.Cells.Item(234, 13).Value = 470
$ws.Cells.Item(234, 14).Value = 3800
$ws.Cells.Item(234, 15).Value = 4000
$ws.Cells.Item(234, 16).Value = 3894
$ws.Cells.Item(234, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(234, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(234, 19).Value = 1947
$ws.Cells.Item(234, 20).Value = 2
